$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

# Delete row 51 entirely - all rows below shift up by one
$ws.Rows.Item(51).Delete()
